# Add a "category" (gender/discipline) column to the partenaire database.
# Column D gets a header "category" in D1 and a category label for each
# data row (rows 2-57), matching the A/B/C row it describes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "category",
    "Environnemental science",
    "Environnemental science",
    "Environnemental science",
    "Environnemental science",
    "Environnemental science",
    "Physics",
    "Chemestry",
    "Physics",
    "Physics",
    "Social Sciences",
    "Biology",
    "Biology",
    "Biology",
    "Medecine",
    "Medecine",
    "Biology",
    "Medecine",
    "Medecine",
    "Medecine",
    "Biology",
    "Biology",
    "Geography",
    "Computational Science",
    "Computational Science",
    "Computational Science",
    "Social Sciences",
    "Social Sciences",
    "Social Sciences",
    "Chemestry",
    "Physics",
    "Physics",
    "Social Sciences",
    "Engineering",
    "Medecine",
    "Medecine",
    "Medecine",
    "Neurosciences",
    "Social Sciences",
    "Engineering",
    "Mathematics",
    "Social Sciences",
    "Physique",
    "Biology",
    "Biology",
    "Mathematics",
    "Mathematics",
    "Physics",
    "Computational Science",
    "Environnemental science",
    "Physics",
    "Engineering",
    "Medecine",
    "Social Sciences",
    "Social Sciences",
    "Social Sciences",
    "Environnemental science"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# New column C got some breathing room when D appeared; widen it to match.
$ws.Columns.Item(3).ColumnWidth = 28.8333333333333

# Keep the selection near where the last edit happened, like the live
# session that produced this change.
$ws.Application.ActiveWindow.ScrollRow = 44
$ws.Range("D9").Select()
